$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localization rows for the "calculate damage" feature.
$ws.Range("A150").Value = "Home_Index_SearchAttack"
$ws.Range("B150").Value = "Attacke suchen"

$ws.Range("A151").Value = "Home_Index_MinDamage"
$ws.Range("B151").Value = "Minimaler Schaden"

$ws.Range("A152").Value = "Home_Index_MaxDamage"
$ws.Range("B152").Value = "Maximaler Schaden"

# Scroll / select to match where the author's cursor ended up after the edit.
$excel.ActiveWindow.ScrollRow = 127
$ws.Range("A151:B152").Select()
